$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6580477.5
$ws.Range("I19").Value = 13889532
$ws.Range("J19").Value = 2327.7
$ws.Range("K19").Value = 13889532
$ws.Range("L19").Value = 2327.7
$ws.Range("M19").Value = -13889357
$ws.Range("N19").Value = -2677.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 28461212
$ws.Range("I98").Value = 10001516
$ws.Range("J98").Value = 102300000
$ws.Range("K98").Value = 10001516
$ws.Range("L98").Value = 102300000
$ws.Range("M98").Value = -10000018
$ws.Range("N98").Value = -102302996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1390755.5
$ws.Range("I112").Value = 1795
$ws.Range("J112").Value = 1604441.8
$ws.Range("K112").Value = 5385
$ws.Range("L112").Value = 4813325.4
$ws.Range("M112").Value = -4277
$ws.Range("N112").Value = -4815541.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5091475
$ws.Range("I116").Value = 2307033
$ws.Range("J116").Value = 10486332
$ws.Range("K116").Value = 2307033
$ws.Range("L116").Value = 10486332
$ws.Range("M116").Value = -2303591
$ws.Range("N116").Value = -10493216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 28461212
$ws.Range("I122").Value = 10001516
$ws.Range("J122").Value = 102300000
$ws.Range("K122").Value = 30004548
$ws.Range("L122").Value = 306900000
$ws.Range("M122").Value = -30002098
$ws.Range("N122").Value = -306904900

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 21045.834
$ws.Range("J130").Value = 21045.834
$ws.Range("L130").Value = 21045.834
$ws.Range("N130").Value = -31085.834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2838811.5
$ws.Range("I132").Value = 679364.9
$ws.Range("J132").Value = 10102405
$ws.Range("K132").Value = 2038094.7
$ws.Range("L132").Value = 30307215
$ws.Range("M132").Value = -2035564.7
$ws.Range("N132").Value = -30312275

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3170.1482
$ws.Range("I141").Value = 1889.9524
$ws.Range("J141").Value = 7650.8335
$ws.Range("K141").Value = 5669.857199999999
$ws.Range("L141").Value = 22952.5005
$ws.Range("M141").Value = -489.8571999999995
$ws.Range("N141").Value = -33312.50049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11028.147
$ws.Range("I2").Value = 14744.2
$ws.Range("J2").Value = 705.7778
$ws.Range("K2").Value = 14744.2
$ws.Range("L2").Value = 705.7778
$ws.Range("M2").Value = -14631.2
$ws.Range("N2").Value = -931.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2513826.2
$ws.Range("I61").Value = 1437952.5
$ws.Range("J61").Value = 5350220.5
$ws.Range("K61").Value = 1437952.5
$ws.Range("L61").Value = 5350220.5
$ws.Range("M61").Value = -1437740.5
$ws.Range("N61").Value = -5350644.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1725.5714
$ws.Range("I63").Value = 1648.1034
$ws.Range("J63").Value = 2100
$ws.Range("K63").Value = 1648.1034
$ws.Range("L63").Value = 2100
$ws.Range("M63").Value = -962.1034
$ws.Range("N63").Value = -3472

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1725.5714
$ws.Range("I66").Value = 1648.1034
$ws.Range("J66").Value = 2100
$ws.Range("K66").Value = 8240.517
$ws.Range("L66").Value = 10500
$ws.Range("M66").Value = -4808.517
$ws.Range("N66").Value = -17364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 645.13635
$ws.Range("I97").Value = 717.1429000000001
$ws.Range("J97").Value = 611.5333000000001
$ws.Range("K97").Value = 717.1429000000001
$ws.Range("L97").Value = 611.5333000000001
$ws.Range("M97").Value = -221.1429000000001
$ws.Range("N97").Value = -1603.5333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1794
$ws.Range("I110").Value = 1103.1818
$ws.Range("J110").Value = 4327
$ws.Range("K110").Value = 1103.1818
$ws.Range("L110").Value = 4327
$ws.Range("M110").Value = 941.8181999999999
$ws.Range("N110").Value = -8417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 11028.147
$ws.Range("I116").Value = 14744.2
$ws.Range("J116").Value = 705.7778
$ws.Range("K116").Value = 14744.2
$ws.Range("L116").Value = 705.7778
$ws.Range("M116").Value = -12450.2
$ws.Range("N116").Value = -5293.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 16576267
$ws.Range("I132").Value = 17927708
$ws.Range("J132").Value = 6947254.5
$ws.Range("K132").Value = 53783124
$ws.Range("L132").Value = 20841763.5
$ws.Range("M132").Value = -53780594
$ws.Range("N132").Value = -20846823.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2513826.2
$ws.Range("I136").Value = 1437952.5
$ws.Range("J136").Value = 5350220.5
$ws.Range("K136").Value = 4313857.5
$ws.Range("L136").Value = 16050661.5
$ws.Range("M136").Value = -4311307.5
$ws.Range("N136").Value = -16055761.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11028.147
$ws.Range("I3").Value = 14744.2
$ws.Range("J3").Value = 705.7778
$ws.Range("K3").Value = 14744.2
$ws.Range("L3").Value = 705.7778
$ws.Range("M3").Value = -14630.2
$ws.Range("N3").Value = -933.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1883.2
$ws.Range("I86").Value = 1915.6875
$ws.Range("J86").Value = 1103.5
$ws.Range("K86").Value = 1915.6875
$ws.Range("L86").Value = 1103.5
$ws.Range("M86").Value = -792.6875
$ws.Range("N86").Value = -3349.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1883.2
$ws.Range("I89").Value = 1915.6875
$ws.Range("J89").Value = 1103.5
$ws.Range("K89").Value = 9578.4375
$ws.Range("L89").Value = 5517.5
$ws.Range("M89").Value = -3962.4375
$ws.Range("N89").Value = -16749.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1053363.9
$ws.Range("I22").Value = 1538664.6
$ws.Range("J22").Value = 1879
$ws.Range("K22").Value = 1538664.6
$ws.Range("L22").Value = 1879
$ws.Range("M22").Value = -1538314.6
$ws.Range("N22").Value = -2579

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 939839.6
$ws.Range("I31").Value = 3088220.8
$ws.Range("J31").Value = 4254.3066
$ws.Range("K31").Value = 3088220.8
$ws.Range("L31").Value = 4254.3066
$ws.Range("M31").Value = -3087925.8
$ws.Range("N31").Value = -4844.3066

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 939839.6
$ws.Range("I34").Value = 3088220.8
$ws.Range("J34").Value = 4254.3066
$ws.Range("K34").Value = 3088220.8
$ws.Range("L34").Value = 4254.3066
$ws.Range("M34").Value = -3088018.8
$ws.Range("N34").Value = -4658.3066

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 21743610
$ws.Range("I94").Value = 1313.75
$ws.Range("J94").Value = 26320934
$ws.Range("K94").Value = 1313.75
$ws.Range("L94").Value = 26320934
$ws.Range("M94").Value = -862.75
$ws.Range("N94").Value = -26321836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 946.9583
$ws.Range("I107").Value = 450.16666
$ws.Range("J107").Value = 1112.5555
$ws.Range("K107").Value = 450.16666
$ws.Range("L107").Value = 1112.5555
$ws.Range("M107").Value = 1469.83334
$ws.Range("N107").Value = -4952.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4119.3716
$ws.Range("I122").Value = 5625.273
$ws.Range("J122").Value = 1570.9231
$ws.Range("K122").Value = 16875.819
$ws.Range("L122").Value = 4712.7693
$ws.Range("M122").Value = -14425.819
$ws.Range("N122").Value = -9612.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1567.1111
$ws.Range("I132").Value = 1067.7838
$ws.Range("J132").Value = 3876.5
$ws.Range("K132").Value = 3203.3514
$ws.Range("L132").Value = 11629.5
$ws.Range("M132").Value = -673.3513999999996
$ws.Range("N132").Value = -16689.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1746558.9
$ws.Range("I134").Value = 8120.4375
$ws.Range("J134").Value = 5720132.5
$ws.Range("K134").Value = 24361.3125
$ws.Range("L134").Value = 17160397.5
$ws.Range("M134").Value = -21826.3125
$ws.Range("N134").Value = -17165467.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1359318.2
$ws.Range("I5").Value = 432.90625
$ws.Range("J5").Value = 4465342
$ws.Range("K5").Value = 1298.71875
$ws.Range("L5").Value = 13396026
$ws.Range("M5").Value = -1186.71875
$ws.Range("N5").Value = -13396250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1359318.2
$ws.Range("I135").Value = 432.90625
$ws.Range("J135").Value = 4465342
$ws.Range("K135").Value = 3896.15625
$ws.Range("L135").Value = 40188078
$ws.Range("M135").Value = -1361.15625
$ws.Range("N135").Value = -40193148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5012506
$ws.Range("I70").Value = 2159079.8
$ws.Range("J70").Value = 11908286
$ws.Range("K70").Value = 2159079.8
$ws.Range("L70").Value = 11908286
$ws.Range("M70").Value = -2158809.8
$ws.Range("N70").Value = -11908826

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5012506
$ws.Range("I73").Value = 2159079.8
$ws.Range("J73").Value = 11908286
$ws.Range("K73").Value = 2159079.8
$ws.Range("L73").Value = 11908286
$ws.Range("M73").Value = -2158143.8
$ws.Range("N73").Value = -11910158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8873.5
$ws.Range("I126").Value = 11241.4
$ws.Range("J126").Value = 2953.75
$ws.Range("K126").Value = 33724.2
$ws.Range("L126").Value = 8861.25
$ws.Range("M126").Value = -31254.2
$ws.Range("N126").Value = -13801.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 50004024
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 71434040
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 71434040
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -71434630

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 50004024
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 71434040
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 71434040
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -71434254

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 71429370
$ws.Range("I46").Value = 1249.5
$ws.Range("J46").Value = 100000620
$ws.Range("K46").Value = 1249.5
$ws.Range("L46").Value = 100000620
$ws.Range("M46").Value = -1061.5
$ws.Range("N46").Value = -100000996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 54980
$ws.Range("J130").Value = 54980
$ws.Range("L130").Value = 54980
$ws.Range("N130").Value = -65020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 15625792
$ws.Range("I126").Value = 27778320
$ws.Range("J126").Value = 1115
$ws.Range("K126").Value = 83334960
$ws.Range("L126").Value = 3345
$ws.Range("M126").Value = -83332490
$ws.Range("N126").Value = -8285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 838530.2
$ws.Range("I132").Value = 3724.25
$ws.Range("J132").Value = 1644549.8
$ws.Range("K132").Value = 11172.75
$ws.Range("L132").Value = 4933649.4
$ws.Range("M132").Value = -8642.75
$ws.Range("N132").Value = -4938709.4
